$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- First, replicate the header/label styling onto the new cells ---
# B1 already carries the bold/bordered/centered style; copy it to C1 and D1.
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial($xlPasteFormats)

# A2 already carries the row-label style; copy it onto the other label cells
# that need it (A3 keeps its own style, A4/A5/A6 are new/shifted rows).
$ws.Range("A2").Copy()
$ws.Range("A4:A6").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# --- Now fill in all the values ---

# Row 1 headers
$ws.Range("B1").Value = "Exp 1"
$ws.Range("C1").Value = "Exp 2"
$ws.Range("D1").Value = "Exp 3"

# Row 2: Accuracy
$ws.Range("A2").Value = "Accuracy"
$ws.Range("B2").Value = 0.9707317073170731
$ws.Range("C2").Value = 0.9609756097560975
$ws.Range("D2").Value = 0.9512195121951219

# Row 3: Error Rate (was Sensitivity)
$ws.Range("A3").Value = "Error Rate"
$ws.Range("B3").Value = 0.02926829268292686
$ws.Range("C3").Value = 0.03902439024390247
$ws.Range("D3").Value = 0.04878048780487809

# Row 4: Sensitivity (new row)
$ws.Range("A4").Value = "Sensitivity"
$ws.Range("B4").Value = 0.9436619718309859
$ws.Range("C4").Value = 0.9354838709677419
$ws.Range("D4").Value = 0.9342105263157895

# Row 5: Specificity (shifted down from row 4)
$ws.Range("A5").Value = "Specificity"
$ws.Range("B5").Value = 0.9850746268656716
$ws.Range("C5").Value = 0.972027972027972
$ws.Range("D5").Value = 0.9612403100775194

# Row 6: Geometric Mean (new row)
$ws.Range("A6").Value = "Geometric Mean"
$ws.Range("B6").Value = 0.9641459769084411
$ws.Range("C6").Value = 0.9535808775146716
$ws.Range("D6").Value = 0.9476290497834435
